$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that should no longer exist ---
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()

# --- Copy formatting for newly introduced cells from a same-column reference cell ---
$ws.Range("B10").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("A10").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A10").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Set changed cell values (only cells whose content actually changes) ---
$ws.Range("B10").Value = "Como parte fundamental da formação específica e geral, a disciplina tem por objetivos (a) fornecer conceitos fundamentais de circuitos elétricos que são importantes para a formação em engenharia física; (b) capacitar o aluno, trabalhando individualmente e em grupo, a modelar e resolver problemas de interesse em circuitos elétricos passivos, com escolha adequada de hipóteses e aplicação de ferramentas correspondentes de solução; (c) introduzir os componentes, técnicas, softwares e equipamentos utilizados na análise e projeto de circuitos elétricos; e (d) aplicar e estender os conceitos físicos aprendidos previamente."
$ws.Range("C10").Value = "Como parte fundamental da formação específica e geral, a disciplina tem por objetivos (a) fornecer conceitos fundamentais de circuitos elétricos que são importantes para a formação em engenharia física; (b) capacitar o aluno, trabalhando individualmente e em grupo, a modelar e resolver problemas de interesse em circuitos elétricos passivos, com escolha adequada de hipóteses e aplicação de ferramentas correspondentes de solução; (c) introduzir os componentes, técnicas, softwares e equipamentos utilizados na análise e projeto de circuitos elétricos; e (d) aplicar e estender os conceitos físicos aprendidos previamente."
$ws.Range("B13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("B14").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C14").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("A15").Value = "Programa resumido:"
$ws.Range("B15").Value = "Introdução à teoria de circuitos elétricos. Elementos de circuitos lineares. Leis fundamentais. Teoremas de circuitos elétricos. Capacitores e indutores. Análise de circuitos DC e AC nos domínios do tempo e da frequência.Identificação de resistores, capacitores e indutores. Medidas de resistência, tensão e corrente elétrica. Utilização de softwares de simulação e projeto de circuitos eletrônicos. Análise e projeto de circuitos DC e AC nos domínios do tempo e da frequência."
$ws.Range("C15").Value = "Introdução à teoria de circuitos elétricos. Elementos de circuitos lineares. Leis fundamentais. Teoremas de circuitos elétricos. Capacitores e indutores. Análise de circuitos DC e AC nos domínios do tempo e da frequência.Identificação de resistores, capacitores e indutores. Medidas de resistência, tensão e corrente elétrica. Utilização de softwares de simulação e projeto de circuitos eletrônicos. Análise e projeto de circuitos DC e AC nos domínios do tempo e da frequência."
$ws.Range("A16").Value = "Short syllabus:"
$ws.Range("B16").Value = "Elements of linear circuits. Ohm's Law. Kirchhoff's Laws. Analysis methods. Theorems. First order circuits. sinusoidal excitation. phasors. Analysis in steady state C.A. Power in steady state A.C. Three Phase Circuits. Transformers. engines. Computational simulation."
$ws.Range("C16").Value = "Elements of linear circuits. Ohm's Law. Kirchhoff's Laws. Analysis methods. Theorems. First order circuits. sinusoidal excitation. phasors. Analysis in steady state C.A. Power in steady state A.C. Three Phase Circuits. Transformers. engines. Computational simulation."
$ws.Range("A17").Value = "Programa:"
$ws.Range("B17").Value = "1. Elementos de Circuitos Lineares; Lei de Ohm; Potência Elétrica.    2. Leis de Kirchhoff; Divisores de Tensão e Corrente; Resistência Série e Paralela; Conversões Y-Delta.    3. Análise Nodal; Análise de Malhas.    4. Teoremas: Linearidade e Superposição, Transformação de Fontes, Thévenin e Norton, Máxima Transferência de Potência.    5. Resistores; Código de Cores; Medições de Tensão e Corrente Elétrica; Verificação Experimental das Leis de Kirchhoff e de Teoremas de Circuitos.    6. Capacitores; Indutores; Circuitos de Primeira Ordem; Circuitos de Segunda Ordem;    7. Simulador SPICE; Análise Computacional de Circuitos de Primeira e Segunda Ordem.    8. Senoides e Fasores; Impedância e Admitância; Leis de Kirchhoff no Domínio da Frequência.    9. Função de Transferência; Ganho e Fase de Tensão/Corrente AC; Ressonância; Filtros.    10. Análise de Potência em CA; Valores Eficazes de Tensão e Corrente; Potência Aparente e Fator de Potência; Potência Complexa; Indutância Mútua; Transformador Ideal.    11. Osciloscópio; Medidas de Amplitude, Frequência e Fase; Análise Experimental de Filtros Passivos."
$ws.Range("C17").Value = "1. Elementos de Circuitos Lineares; Lei de Ohm; Potência Elétrica.    2. Leis de Kirchhoff; Divisores de Tensão e Corrente; Resistência Série e Paralela; Conversões Y-Delta.    3. Análise Nodal; Análise de Malhas.    4. Teoremas: Linearidade e Superposição, Transformação de Fontes, Thévenin e Norton, Máxima Transferência de Potência.    5. Resistores; Código de Cores; Medições de Tensão e Corrente Elétrica; Verificação Experimental das Leis de Kirchhoff e de Teoremas de Circuitos.    6. Capacitores; Indutores; Circuitos de Primeira Ordem; Circuitos de Segunda Ordem;    7. Simulador SPICE; Análise Computacional de Circuitos de Primeira e Segunda Ordem.    8. Senoides e Fasores; Impedância e Admitância; Leis de Kirchhoff no Domínio da Frequência.    9. Função de Transferência; Ganho e Fase de Tensão/Corrente AC; Ressonância; Filtros.    10. Análise de Potência em CA; Valores Eficazes de Tensão e Corrente; Potência Aparente e Fator de Potência; Potência Complexa; Indutância Mútua; Transformador Ideal.    11. Osciloscópio; Medidas de Amplitude, Frequência e Fase; Análise Experimental de Filtros Passivos."
$ws.Range("A18").Value = "Syllabus:"
$ws.Range("B18").Value = "1. Linear Circuit Elements; Ohm's Law; Electric power.2. Kirchhoff's Laws; Voltage and Current Networks; Series and Parallel Resistance; Y-Delta Conversions.3. Nodal Analysis; Mesh Analysis.4. Theorems: Linearity and Superposition, Source Transformation, Thévenin and Norton, Maximum Power Transfer.5. Resistors; Color Code; Voltage and Electric Current Measurements; Experimental Verification of Kirchhoff's Laws and Circuit Theorems.6. Capacitors; Inductors; First Order Circuits; Second Order Circuits;7. SPICE simulator; Computational Analysis of First and Second Order Circuits.8. Sinusoids and Phasors; Impedance and Admittance; Kirchhoff's Laws in the Frequency Domain.9. Transfer Function; AC Voltage/Current Gain and Phase; Resonance; filters.10. AC Power Analysis; Effective Voltage and Current Values; Apparent Power and Power Factor; Complex Power; Mutual Inductance; Ideal Transformer.11. Oscilloscope; Amplitude, Frequency and Phase Measurements; Experimental Analysis of Passive Filters."
$ws.Range("C18").Value = "1. Linear Circuit Elements; Ohm's Law; Electric power.2. Kirchhoff's Laws; Voltage and Current Networks; Series and Parallel Resistance; Y-Delta Conversions.3. Nodal Analysis; Mesh Analysis.4. Theorems: Linearity and Superposition, Source Transformation, Thévenin and Norton, Maximum Power Transfer.5. Resistors; Color Code; Voltage and Electric Current Measurements; Experimental Verification of Kirchhoff's Laws and Circuit Theorems.6. Capacitors; Inductors; First Order Circuits; Second Order Circuits;7. SPICE simulator; Computational Analysis of First and Second Order Circuits.8. Sinusoids and Phasors; Impedance and Admittance; Kirchhoff's Laws in the Frequency Domain.9. Transfer Function; AC Voltage/Current Gain and Phase; Resonance; filters.10. AC Power Analysis; Effective Voltage and Current Values; Apparent Power and Power Factor; Complex Power; Mutual Inductance; Ideal Transformer.11. Oscilloscope; Amplitude, Frequency and Phase Measurements; Experimental Analysis of Passive Filters."
$ws.Range("A19").Value = "Avaliação:"
$ws.Range("A20").Value = "Método:"
$ws.Range("B20").Value = "Aulas expositivas e práticas de laboratório com interações em grupo para a solução de problemas."
$ws.Range("C20").Value = "Aulas expositivas e práticas de laboratório com interações em grupo para a solução de problemas."
$ws.Range("A21").Value = "Critério:"
$ws.Range("B21").Value = "Média aritmética (M) de provas individuais (P1 e P2) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,3*P1+0,3*P2+0.4*T"
$ws.Range("C21").Value = "Média aritmética (M) de provas individuais (P1 e P2) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,3*P1+0,3*P2+0.4*T"
$ws.Range("A22").Value = "Norma de recuperação:"
$ws.Range("B22").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C22").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("A23").Value = "Bibliografia:"
$ws.Range("B23").Value = "JOHNSON, D. E. Fundamentos de análise de circuitos elétricos. Prentice Hall, 1994.HAYT, W. H. Análise de circuitos em engenharia. McGraw-Hill, 2008.DESOER, C. A. Teoria básica de circuitos. Guanabara Dois, 1979.SCOTT, R. E. Elements of linear circuits. Addison Wesley, 1965ALEXANDER, C. K. E SADIKU, M. N. O. Fundamentos de Circuitos Elétricos. McGraw-Hill, 2013.NILSSON, J. W. E RIEDEL, S. A. Electric Circuits. Prentice Hall, 2011.BOYLESTAD, Robert L. Introdução à Análise de Circuitos. Pearson, 2011"
$ws.Range("C23").Value = "JOHNSON, D. E. Fundamentos de análise de circuitos elétricos. Prentice Hall, 1994.HAYT, W. H. Análise de circuitos em engenharia. McGraw-Hill, 2008.DESOER, C. A. Teoria básica de circuitos. Guanabara Dois, 1979.SCOTT, R. E. Elements of linear circuits. Addison Wesley, 1965ALEXANDER, C. K. E SADIKU, M. N. O. Fundamentos de Circuitos Elétricos. McGraw-Hill, 2013.NILSSON, J. W. E RIEDEL, S. A. Electric Circuits. Prentice Hall, 2011.BOYLESTAD, Robert L. Introdução à Análise de Circuitos. Pearson, 2011"
$ws.Range("A24").Value = "Requisitos:"
$ws.Range("B25").Value = "LOB1006 -  Cálculo IV  (Requisito)`n"
$ws.Range("C25").Value = "LOB1006 -  Cálculo IV  (Requisito)`n"

# --- Set row heights ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(25).RowHeight = 30
